# Adds a "Frase traducida" (translation) column D to the WriteActivity
# exercise sheet, plus extra rows (8-14) so every one of the 13 existing
# phrase/word rows gets its own translation placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 8-14: repeat the Frases/Palabras pattern of rows 6-7 ----
# (row 6 pattern -> even "Smart money..." entries, row 7 pattern -> odd
#  "!chan chan chaaann!" entries), just bumping the level number in A.
$rowSixB  = $ws.Range("B6").Value2
$rowSixC  = $ws.Range("C6").Value2
$rowSevenB = $ws.Range("B7").Value2
$rowSevenC = $ws.Range("C7").Value2

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = $rowSixB
$ws.Range("C8").Value = $rowSixC

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = $rowSevenB
$ws.Range("C9").Value = $rowSevenC

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = $rowSixB
$ws.Range("C10").Value = $rowSixC

$ws.Range("A11").Value = 7
$ws.Range("B11").Value = $rowSevenB
$ws.Range("C11").Value = $rowSevenC

$ws.Range("A12").Value = 8
$ws.Range("B12").Value = $rowSixB
$ws.Range("C12").Value = $rowSixC

$ws.Range("A13").Value = 9
$ws.Range("B13").Value = $rowSixB
$ws.Range("C13").Value = $rowSixC

$ws.Range("A14").Value = 10
$ws.Range("B14").Value = $rowSevenB
$ws.Range("C14").Value = $rowSevenC

# --- Formatting: copy the row-6 / row-7 styles down onto the new rows -
$ws.Range("A6:C6").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A13:C13").PasteSpecial(-4122)

$ws.Range("A7:C7").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A14:C14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- New column D: header + one translation placeholder per row -------
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D1").Value = "Frase traducida"

$ws.Range("D2").Value = "Traduccion 1"
$ws.Range("D3").Value = "Traduccion 2"
$ws.Range("D4").Value = "Traduccion 3"
$ws.Range("D5").Value = "Traduccion 4"
$ws.Range("D6").Value = "Traduccion 5"
$ws.Range("D7").Value = "Traduccion 6"
$ws.Range("D8").Value = "Traduccion 7"
$ws.Range("D9").Value = "Traduccion 8"
$ws.Range("D10").Value = "Traduccion 9"
$ws.Range("D11").Value = "Traduccion 10"
$ws.Range("D12").Value = "Traduccion 11"
$ws.Range("D13").Value = "Traduccion 12"
$ws.Range("D14").Value = "Traduccion 13"

# --- Column width for the new column, matching the others -------------
$ws.Columns.Item(4).ColumnWidth = 16.140625

# --- Sheet view / selection, per the workbook author's last action ----
$ws.Range("D21").Select()

Write-Output "ok"
